$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '248.41'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.55'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05683'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.411'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.325'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8072'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8986'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1418'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07442'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03055'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'ProBitToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1342'
$ws.Range("E13").Value = '12ProBitTokenPROB'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03032'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09394'
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.883'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001584'
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.04784'
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("B19").Value = 'UpBots'
$ws.Range("C19").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.01829'
$ws.Range("E19").Value = '18UpBotsUBXTBestin24h'
$ws.Range("B20").Value = 'One'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0005804'
$ws.Range("E20").Value = '19OneONE'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.006425'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.004985'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0009969'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0001500'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.170'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03972'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.003046'
$ws.Range("E41").Value = '40KickTokenKICKWorstin24h'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1071'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002731'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007695'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005588'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000751'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4994'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2027'
$ws.Range("E48").Value = '47BOLOBOLO'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002102'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.01011'
